$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.411943666666667
$ws.Range("H2").Value = 4.235831
$ws.Range("M2").Value = 18.28424166666667
$ws.Range("N2").Value = 54.85272500000001
$ws.Range("O2").Value = 0.9545246501532072
$ws.Range("P2").Value = 0.9545246501532071
$ws.Range("Q2").Value = 25.81631922105278
$ws.Range("R2").Value = 232.346872989475
$ws.Range("S2").Value = 0.9545246501532072
$ws.Range("T2").Value = 0.9545246501532071

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.411943666666667
$ws.Range("H3").Value = 4.235831
$ws.Range("O3").Value = 0.02163433578790886
$ws.Range("P3").Value = 0.02163433578790886
$ws.Range("Q3").Value = 0.5851278107343333
$ws.Range("R3").Value = 5.266150296609
$ws.Range("S3").Value = 0.02163433578790886
$ws.Range("T3").Value = 0.02163433578790886

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.411943666666667
$ws.Range("H4").Value = 4.235831
$ws.Range("O4").Value = 0.02384101405888406
$ws.Range("P4").Value = 0.02384101405888406
$ws.Range("Q4").Value = 0.6448101988764445
$ws.Range("R4").Value = 5.803291789888001
$ws.Range("S4").Value = 0.02384101405888406
$ws.Range("T4").Value = 0.02384101405888406
